$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A19").Value = "Dear Sir"
$ws.Range("A19").Select() | Out-Null
